$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q3" worksheet right after "总计", by copying
#    the existing "2022-Q2" sheet so it inherits identical formatting,
#    then renaming it and overwriting its data with the 2022-Q3 figures.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($null, $totalSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Force text storage for the text-like columns (fund code / figures are
# stored as text in this workbook), then strip the temporary number
# format back off so no stray style is left behind on the cells.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "010695"
$newSheet.Range("C2").Value = "华夏磐益一年定期开放混合"
$newSheet.Range("D2").Value = "15.90"
$newSheet.Range("E2").Value = "99.95"
$newSheet.Range("F2").Value = "4.21"
$newSheet.Range("G2").Value = "0.6694"
$newSheet.Range("H2").Value = 7

$newSheet.Range("H2").Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)
$newSheet.Range("A1").Select() | Out-Null

# ------------------------------------------------------------------
# 2) Update the "总计" summary sheet: add a new top data row for
#    2022-Q3 and push the existing rows down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Extend the existing formatting down into the new row 6 first.
$total.Range("A5:D5").Copy()
$total.Range("A6:D6").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.67

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 2.77

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.11

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 9
$total.Range("D5").Value = 1.3

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.97

$total.Range("A1").Select() | Out-Null
